$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A2:A41 with the single text value "28-04-2025" (stored as text, not a date)
$rng = $ws.Range("A2:A41")
$rng.Value = "28-04-2025"

# Update selection to match the post-edit state
$ws.Range("A4:A41").Select()
